# Applies the crypto-price/volume refresh described by the commit.
# Cells whose new text looks like a plain number (e.g. "1.00", "7.63")
# are written with a leading apostrophe so Excel keeps them as TEXT
# (matching the original inline-string cells) instead of silently
# converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.997.56'
$ws.Range('E2').Value = '  +1.01%  '

# Row 3
$ws.Range('D3').Value = '3.379.80'

# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').Value = '''570.01'
$ws.Range('E5').Value = '  +0.50%  '

# Row 6
$ws.Range('D6').Value = '''140.36'
$ws.Range('E6').Value = '  +0.33%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('E8').Value = '  +0.19%  '

# Row 9
$ws.Range('D9').Value = '''7.63'
$ws.Range('E9').Value = '  +1.58%  '

# Row 10
$ws.Range('E10').Value = '  -1.16%  '

# Row 11
$ws.Range('E11').Value = '  -0.23%  '

# Row 12
$ws.Range('D12').Value = '3.957.43'
$ws.Range('E12').Value = '  +0.16%  '

# Row 13
$ws.Range('E13').Value = '  +1.98%  '

# Row 14
$ws.Range('D14').Value = '''27.82'
$ws.Range('E14').Value = '  -0.74%  '

# Row 15
$ws.Range('D15').Value = '3.368.96'
$ws.Range('E15').Value = '  -0.38%  '

# Row 16
$ws.Range('E16').Value = '  -0.33%  '

# Row 17
$ws.Range('D17').Value = '61.100.29'
$ws.Range('E17').Value = '  +0.93%  '

# Row 18
$ws.Range('E18').Value = '  -1.47%  '

# Row 19
$ws.Range('D19').Value = '''13.52'
$ws.Range('E19').Value = '  -1.78%  '

# Row 20
$ws.Range('D20').Value = '''8.88'
$ws.Range('E20').Value = '  -1.20%  '

# Row 21
$ws.Range('D21').Value = '''381.98'
$ws.Range('E21').Value = '  -1.06%  '

# Row 22
$ws.Range('E22').Value = '  +3.26%  '

# Row 23
$ws.Range('D23').Value = '''0.549'
$ws.Range('E23').Value = '  -0.80%  '

# Row 24
$ws.Range('E24').Value = '  -0.19%  '

# Row 25
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.515.66'
$ws.Range('E25').Value = '  -0.07%  '

# Row 26
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '''0.0000113'
$ws.Range('E26').Value = '  -1.77%  '

# Row 27
$ws.Range('D27').Value = '''0.193'
$ws.Range('E27').Value = '  +8.22%  '

# Row 28
$ws.Range('D28').Value = '''1.01'
$ws.Range('E28').Value = '  +1.00%  '

# Row 29
$ws.Range('D29').Value = '''7.19'
$ws.Range('E29').Value = '  -2.20%  '

# Row 30
$ws.Range('D30').Value = '''7.94'
$ws.Range('E30').Value = '  +0.52%  '

# Row 31
$ws.Range('E31').Value = '  -0.23%  '

# Row 32
$ws.Range('E32').Value = '  -0.04%  '

# Row 33
$ws.Range('E33').Value = '  -3.64%  '

# Row 34
$ws.Range('D34').Value = '''23.18'
$ws.Range('E34').Value = '  -1.14%  '

# Row 35
$ws.Range('E35').Value = '  +0.73%  '

# Row 36
$ws.Range('D36').Value = '''166.22'
$ws.Range('E36').Value = '  -1.26%  '

# Row 37
$ws.Range('D37').Value = '3.415.10'
$ws.Range('E37').Value = '  +0.29%  '

# Row 38
$ws.Range('E38').Value = '  +0.78%  '

# Row 39
$ws.Range('D39').Value = '''1.46'
$ws.Range('E39').Value = '  -1.90%  '

# Row 40
$ws.Range('D40').Value = '''0.0763'
$ws.Range('E40').Value = '  -0.60%  '

# Row 41
$ws.Range('D41').Value = '''25.83'
$ws.Range('E41').Value = '  -4.06%  '

# Row 42
$ws.Range('E42').Value = '  +0.05%  '

# Row 43
$ws.Range('E43').Value = '  +0.22%  '

# Row 44
$ws.Range('E44').Value = '  -1.58%  '

# Row 45
$ws.Range('E45').Value = '  -3.01%  '

# Row 46
$ws.Range('E46').Value = '  -0.58%  '

# Row 47
$ws.Range('D47').Value = '2.437.10'
$ws.Range('E47').Value = '  -3.10%  '

# Row 48
$ws.Range('D48').Value = '''22.77'
$ws.Range('E48').Value = '  -1.27%  '

# Row 49
$ws.Range('E49').Value = '  -1.77%  '

# Row 50
$ws.Range('D50').Value = '''0.0260'
$ws.Range('E50').Value = '  -2.54%  '

# Row 51
$ws.Range('D51').Value = '''2.08'
$ws.Range('E51').Value = '  +5.15%  '
